# Updated user stories SOI-3792 and SOI-2206.
# Add new rows to the PhoneLine sheet describing the SOI_3792 configuration rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PhoneLine")
$ws.Activate()

$ws.Range("B4").Value = "New,NotApplicable,898989,ADSL+E,NotApplicable,NotApplicable,NotApplicable,Classic Telephone Line"
$ws.Range("A5").Value = "SOI_3792"
$ws.Range("B5").Value = "New,NotApplicable,898989,VDSL2,NotApplicable,NotApplicable,NotApplicable,Classic Telephone Line"
$ws.Range("A6").Value = "SOI_3792"
$ws.Range("B6").Value = "New,NotApplicable,898989,Shared Fiber - GPON,NotApplicable,NotApplicable,NotApplicable,IP Telephone Line"
$ws.Range("A4").Value = "SOI_3792_Configuration_1"

$ws.Range("A7").Select()
